# "Corrigiendo problemas con home": the Home/logo group on slide 1 was
# off-center and non-square. Resize + recenter it into a perfect square
# centered on the (square) slide.
#
# Target geometry (from the authored OOXML):
#   <a:off x="1161256" y="1161256"/><a:ext cx="4876800" cy="4876800"/>
# EMU -> points (1 pt = 12700 EMU):
#   off  = 1161256 / 12700 = 91.43748031496062 pt
#   ext  = 4876800 / 12700 = 384.0 pt
#
# The COM layer here rounds Left/Top to the nearest EMU from the point
# value it is given; 91.43749 is the smallest literal that round-trips
# to exactly 1161256 EMU (91.43748031496062 itself rounds down to
# 1161255 because of float truncation along the way).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)   # "Group 33" - the Home logo group

$sh.Left   = 91.43749
$sh.Top    = 91.43749
$sh.Width  = 384.0
$sh.Height = 384.0
